# Update ObjTables / SBtab metadata header strings across all worksheets.
# Adds schema='SBtab', reorders tableFormat='row' right after type='Data',
# and bumps the date stamp from 2020-03-09 .. to 2020-03-10 00:00:34.

$wb = $excel.ActiveWorkbook

$newDate = "2020-03-10 00:00:34"

# Sheet name -> (id, name, hasDocument)
$sheetInfo = @(
    @{ Sheet = "!!Compartment";             Id = "Compartment";             Name = "Compartment";             Doc = "TestModel layout" },
    @{ Sheet = "!!Compound";                Id = "Compound";                Name = "Compound";                Doc = "TestModel layout" },
    @{ Sheet = "!!Definition";              Id = "Definition";              Name = "Definition";              Doc = $null },
    @{ Sheet = "!!Enzyme";                  Id = "Enzyme";                  Name = "Enzyme";                  Doc = $null },
    @{ Sheet = "!!FbcObjective";            Id = "FbcObjective";            Name = "FbcObjective";            Doc = $null },
    @{ Sheet = "!!Gene";                    Id = "Gene";                    Name = "Gene";                    Doc = $null },
    @{ Sheet = "!!Layout";                  Id = "Layout";                  Name = "Layout";                  Doc = "TestModel layout" },
    @{ Sheet = "!!Measurement";             Id = "Measurement";             Name = "Measurement";             Doc = $null },
    @{ Sheet = "!!PbConfig";                Id = "PbConfig";                Name = "PbConfig";                Doc = $null },
    @{ Sheet = "!!Position";                Id = "Position";                Name = "Position";                Doc = $null },
    @{ Sheet = "!!Protein";                 Id = "Protein";                 Name = "Protein";                 Doc = $null },
    @{ Sheet = "!!Quantity";                Id = "Quantity";                Name = "Quantity";                Doc = $null },
    @{ Sheet = "!!QuantityInfo";            Id = "QuantityInfo";            Name = "QuantityInfo";            Doc = $null },
    @{ Sheet = "!!QuantityMatrix";          Id = "QuantityMatrix";          Name = "QuantityMatrix";          Doc = $null },
    @{ Sheet = "!!Reaction";                Id = "Reaction";                Name = "Reaction";                Doc = "TestModel layout" },
    @{ Sheet = "!!ReactionStoichiometry";   Id = "ReactionStoichiometry";   Name = "ReactionStoichiometry";   Doc = $null },
    @{ Sheet = "!!Regulator";               Id = "Regulator";               Name = "Regulator";               Doc = $null },
    @{ Sheet = "!!Relation";                Id = "Relation";                Name = "Relation";                Doc = $null },
    @{ Sheet = "!!Relationship";            Id = "Relationship";            Name = "Relationship";            Doc = $null },
    @{ Sheet = "!!SparseMatrix";            Id = "SparseMatrix";            Name = "SparseMatrix";            Doc = $null },
    @{ Sheet = "!!SparseMatrixColumn";      Id = "SparseMatrixColumn";      Name = "SparseMatrixColumn";      Doc = $null },
    @{ Sheet = "!!SparseMatrixOrdered";     Id = "SparseMatrixOrdered";     Name = "SparseMatrixOrdered";     Doc = $null },
    @{ Sheet = "!!SparseMatrixRow";         Id = "SparseMatrixRow";         Name = "SparseMatrixRow";         Doc = $null },
    @{ Sheet = "!!StoichiometricMatrix";    Id = "StoichiometricMatrix";    Name = "StoichiometricMatrix";    Doc = $null },
    @{ Sheet = "!!rxnconContingencyList";   Id = "rxnconContingencyList";   Name = "rxnconContingencyList";   Doc = $null },
    @{ Sheet = "!!rxnconReactionList";      Id = "rxnconReactionList";      Name = "rxnconReactionList";      Doc = $null }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Sheet)
    $ws.Unprotect()

    # The very first worksheet additionally carries the top-level
    # "!!!ObjTables" banner in A1, which pushes its own
    # "!!ObjTables ... id='Compartment' ..." line down to A2.
    if ($info.Sheet -eq "!!Compartment") {
        $ws.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.8' date='$newDate'"
    }

    $text = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='$($info.Id)' name='$($info.Name)' date='$newDate' objTablesVersion='0.0.8'"
    if ($info.Doc) {
        $text = "$text document='$($info.Doc)'"
    }

    if ($info.Sheet -eq "!!Compartment") {
        $ws.Range("A2").Value = $text
    } else {
        $ws.Range("A1").Value = $text
    }

    $ws.Protect()
}
